$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price (D) and Volume(1h) (E) values in this sheet are stored as plain text
# (inline strings), so force text number format before assigning to prevent Excel
# from auto-converting numeric-looking strings (e.g. "598.11", "0.0000269") into
# real numbers, which would drop significant digits/trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.821.81"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.195.02"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.11"
$ws.Range("E5").Value = "  +3.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.68"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.188.47"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.06"
$ws.Range("E11").Value = "  -2.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.510"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.36"
$ws.Range("E14").Value = "  +4.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.720.91"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.46"
$ws.Range("E16").Value = "  +3.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.966.66"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.200.42"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.71"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.36"
$ws.Range("E21").Value = "  +3.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.737"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.11"
$ws.Range("E23").Value = "  +3.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.29"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.69"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.00"
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  +3.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.84"
$ws.Range("E30").Value = "  +7.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.86"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.99"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.54"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.97"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0903"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "484.18"
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0418"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("E40").Value = "  -4.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.86"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("E42").Value = "  +5.37%  "
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.946.46"
$ws.Range("E44").Value = "  -3.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0641"
$ws.Range("E45").Value = "  +6.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.37"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("E51").Value = "  +3.54%  "
